# Applies the "StructureDefinition-claim-status" edit:
# adds a new (hidden) "valueCodeableConcept" slice row to the Elements sheet,
# turns the Extension.value[x] row (row 6) into a sliced/closed element,
# and keeps the surrounding ranges (defined name, autofilter, conditional
# formatting, dimension) consistent with the extra row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Edit existing row 6 (Extension.value[x]) to describe the new slicing
# ---------------------------------------------------------------------
$ws.Range("AA6").Value2 = "type:`$this}`n"   # Slicing Discriminator
$ws.Range("AB6").Value2 = ""                  # Slicing Description (blank)
$ws.Range("AD6").Value2 = "closed"            # Slicing Rules

# Remove any auto row-height Excel may have computed because of the
# embedded newline above - the row keeps its default height.
$ws.Rows.Item(6).AutoFit()

# ---------------------------------------------------------------------
# 2. Create row 7 (valueCodeableConcept slice) below it, copying the
#    formatting (style) of row 6 first, then filling in the values.
# ---------------------------------------------------------------------
$ws.Range("A6:AJ6").Copy()
$ws.Range("A7:AJ7").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A7").Value2  = "Extension.value[x]"
$ws.Range("B7").Value2  = "valueCodeableConcept"
$ws.Range("C7").Value2  = ""
$ws.Range("D7").Value2  = ""
$ws.Range("E7").Value2  = "0"
$ws.Range("F7").Value2  = "1"
$ws.Range("G7").Value2  = ""
$ws.Range("H7").Value2  = ""
$ws.Range("I7").Value2  = ""
$ws.Range("J7").Value2  = "CodeableConcept`n"
$ws.Range("K7").Value2  = "Value of extension"
$ws.Range("L7").Value2  = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
$ws.Range("M7").Value2  = ""
$ws.Range("N7").Value2  = ""
$ws.Range("O7").Value2  = ""
$ws.Range("P7").Value2  = ""
$ws.Range("Q7").Value2  = ""
$ws.Range("R7").Value2  = ""
$ws.Range("S7").Value2  = ""
$ws.Range("T7").Value2  = ""
$ws.Range("U7").Value2  = ""
$ws.Range("V7").Value2  = ""
$ws.Range("W7").Value2  = "extensible"
$ws.Range("X7").Value2  = ""
$ws.Range("Y7").Value2  = "https://x12.org/codes/claim-status-codes"
$ws.Range("Z7").Value2  = ""
$ws.Range("AA7").Value2 = ""
$ws.Range("AB7").Value2 = ""
$ws.Range("AC7").Value2 = ""
$ws.Range("AD7").Value2 = ""
$ws.Range("AE7").Value2 = "Extension.value[x]"
$ws.Range("AF7").Value2 = "0"
$ws.Range("AG7").Value2 = "1"
$ws.Range("AH7").Value2 = ""
$ws.Range("AI7").Value2 = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`n"
$ws.Range("AJ7").Value2 = "N/A"

# Remove any auto row-height Excel may have computed because of the
# embedded newline above - the row keeps its default height.
$ws.Rows.Item(7).AutoFit()

# ---------------------------------------------------------------------
# 3. Restore the hidden state of the detail rows (2-7): the import loses
#    the "hidden" flag, so it is re-applied here.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 7; $r++) {
    $ws.Rows.Item($r).Hidden = $true
}

# ---------------------------------------------------------------------
# 4. Resize columns B (Slice Name) and Y (Binding Value Set) to fit the
#    new, longer content.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 21.5
$ws.Columns.Item(25).ColumnWidth = 37.666666666666664

# ---------------------------------------------------------------------
# 5. Extend the AutoFilter range down to the new row, keeping the same
#    per-column filter criteria that were already defined.
# ---------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:AJ7").AutoFilter(7, "<> ", 1)
$ws.Range("A1:AJ7").AutoFilter(27, @(""), 7)

# ---------------------------------------------------------------------
# 6. Extend the conditional formatting range to include the new row.
# ---------------------------------------------------------------------
$fcs = $ws.Range("A2:AI5").FormatConditions
$fcs.Item(1).ModifyAppliesToRange($ws.Range("A2:AI6"))

# ---------------------------------------------------------------------
# 7. Update the workbook-level _FilterDatabase defined name to match the
#    new data extent.
# ---------------------------------------------------------------------
$wb.Names.Item(1).RefersTo = "=Elements!`$A`$1:`$AJ`$7"

Write-Host "edit complete"
